$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).Delete()
$tbl = $s.Shapes.AddTable(8, 2, 23.138346456692915, 18.85677165354331, 913.7233070866142, 502.28645669291336)
$table = $tbl.Table
$cell = $table.Cell(1,1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "Difference between window port and viewport"
$tr.Font.Name = "Times New Roman"
$tr.Font.Size = 32
$tr.Font.Bold = -1
try {
  $tr.Font.Color.RGB = 0
  Write-Host "color set ok"
} catch {
  Write-Host "ERR color: $_"
}
try {
  $tr.ParagraphFormat.Alignment = 2
  Write-Host "align ok"
} catch {
  Write-Host "ERR align: $_"
}
